# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2..16) the existing error values in columns B..(lastCol)
# are shifted one column to the right (B->C, C->D, ... ), and a newly
# computed (near-zero) "Q0" naive-forecast error value is written into the
# now-empty column B. The previous right-most value (old column K, where it
# existed) is dropped, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
$lastCols = @(11,11,11,11,11,10,9,8,7,6,5,4,3,2,1)
$newB = @("0.0000003593882045849206","-0.0000001035781544145298","-0.000000000387512216759589","-0.000000107388789361007","-0.0000001035472805832605","0.000006303355340908645","-0.0000002375649628613696","0.0000003720025918141356","0.0000003829984367986761","-0.000003160475492397508","-0.00000004101096154340844","-0.0000001831659499074156","0.0000002770877186031306","0.000000229775004800814","-0.0000001554241066958895")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $lastCol = $lastCols[$i]

    # The table is bounded at column K (11): any value that would be
    # pushed past column K when shifting right is simply dropped, so the
    # shift only needs to start at column J (10) at most.
    $startCol = $lastCol
    if ($startCol -gt 10) {
        $startCol = 10
    }

    # Shift existing values one column to the right, starting from the
    # right-most column to move and working back down to column B (2),
    # so that each source cell is read before it gets overwritten.
    for ($col = $startCol; $col -ge 2; $col--) {
        $ws.Cells.Item($r, $col + 1).Value = $ws.Cells.Item($r, $col).Value()
    }

    # Write the newly computed value into column B.
    $ws.Cells.Item($r, 2).Value = $newB[$i]
}
